# Register test data workbook update:
#  - add a new "DuplicateEmail" worksheet (after Sheet1/Sheet2) that
#    demonstrates a duplicate-email registration scenario
#  - refresh the saved cursor/selection state on all three sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- New sheet: DuplicateEmail -------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "DuplicateEmail"

$ws3.Range("A1").Value = "FirstName"
$ws3.Range("B1").Value = "LastName"
$ws3.Range("C1").Value = "Email"
$ws3.Range("D1").Value = "Password"
$ws3.Range("E1").Value = "ConfirmPassword"

$ws3.Range("A2").Value = "Ravi"
$ws3.Range("B2").Value = "Kumar"
$ws3.Range("C2").Value = "ravi.kumar1@testmail.com"
$ws3.Range("D2").Value = "Ravi@2024"
$ws3.Range("E2").Value = "Ravi@2024"

$ws3.Range("A3").Value = "Ravi"
$ws3.Range("B3").Value = "Kumar"
$ws3.Range("C3").Value = "ravi.kumar1@testmail.com"
$ws3.Range("D3").Value = "Ravi@2024"
$ws3.Range("E3").Value = "Ravi@2024"

$ws3.Hyperlinks.Add($ws3.Range("C2"), "mailto:ravi.kumar1@testmail.com") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "mailto:ravi.kumar1@testmail.com") | Out-Null

# --- Refresh the selection/cursor state on every sheet --------------------
$ws1.Activate()
$ws1.Range("A1:E3").Select() | Out-Null

$ws2.Activate()
$ws2.Range("A1:E2").Select() | Out-Null

$ws3.Activate()
$ws3.Range("H3").Select() | Out-Null
